$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend formatting for new column W (copy from column V) and new rows 8-11 (copy from row 7) ---
$ws.Range("V1:V7").Copy() | Out-Null
$ws.Range("W1:W7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:W7").Copy() | Out-Null
$ws.Range("A8:W8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:W9").PasteSpecial(-4122) | Out-Null
$ws.Range("A10:W10").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:W11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 1 header numbers (B1:W1) ---
$row1 = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21)
for ($i = 0; $i -lt $row1.Length; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $row1[$i]
}

# --- Row 2: A2 stays 0, B2 stays "HKL"; C2:W2 get HKL-pair labels ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "HKL"
$row2 = @('[1, 1, 1]', '[2, 0, 0]', '[2, 2, 0]', '[3, 1, 1]', '[2, 2, 2]', '[4, 0, 0]', '[3, 3, 1]', '[4, 2, 0]', '[4, 2, 2]', '[5, 1, 1]', '[3, 3, 3]', '1Pair-A', '1Pair-B', '2Pairs-A', '2Pairs-B', '3Pairs-A', '3Pairs-B', '3Pairs-C', '4Pairs', '5A4F', 'MaxUnique')
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 3 + $i).Value = $row2[$i]
}

# --- Column A (rows 3-11): method index 1-9 ---
$colA = @(1, 2, 3, 4, 5, 6, 7, 8, 9)
for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item(3 + $i, 1).Value = $colA[$i]
}

# --- Column B (rows 3-11): method names ---
$colB = @('Equal Angle', 'CLR', 'BT8Hex', 'Spiral', 'OffsetF', 'OffsetA', 'RD Single', 'TD Single', 'HexGrid-90degTilt5degRes')
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item(3 + $i, 2).Value = $colB[$i]
}

# --- Data grid rows 3-11, columns C:W (intensity ratios) ---
$data = @(
    @(1.082161383285303, 0.9005979827089338, 0.9817939481268012, 0.983321325648415, 1.082161383285303, 0.9005979827089338, 1.016087896253602, 0.9638832853025937, 1.034272334293948, 0.9404610951008645, 1.082154178674352, 1.082161383285303, 0.9817939481268012, 0.9411959654178674, 0.9825576368876081, 0.9881844380403457, 0.9552377521613833, 0.9881844380403458, 0.9869686599423632, 1.006007204610951, 0.9878224063400577),
    @(1.001019739608322, 0.9831326726267076, 0.9960509519842181, 0.9920574798733447, 1.001019739608322, 0.9831326726267076, 0.9972227664010135, 0.9939189858118702, 0.9971927827177304, 0.9869878394325744, 1.001027041099043, 1.001019739608322, 0.9960509519842181, 0.9895918123054628, 0.9940542159287814, 0.9934011214064159, 0.9904137014947567, 0.9934011214064159, 0.993065211023148, 0.9946561167401828, 0.9934479023069727),
    @(1.001051086494181, 0.9705521418085087, 1.000695743726981, 0.9913908517389517, 1.001051086494181, 0.9705521418085087, 0.9999661872960939, 0.9945407524359218, 0.9983385337921159, 0.9814411766989068, 1.001042248944182, 1.001051086494181, 1.000695743726981, 0.9856239427677447, 0.9960432977329663, 0.9907663240098902, 0.9875462457581471, 0.9907663240098902, 0.9909224559421556, 0.9929481820525605, 0.9922470592489576),
    @(0.9981523727028009, 0.9928956756264826, 0.9927967605002801, 0.9934011842924885, 0.9981523727028009, 0.9928956756264826, 0.9943708507230946, 0.9939596431722725, 0.9960668685529072, 0.9922132869407635, 0.9981573682987914, 0.9981523727028009, 0.9927967605002801, 0.9928462180633814, 0.9930989723963843, 0.9946149362765212, 0.9930312068064171, 0.9946149362765212, 0.994311498280513, 0.9950796731649707, 0.9942320803138862),
    @(0.9139503744091894, 1.607818008941707, 0.7423551603023174, 1.060944657175708, 0.9139503744091894, 1.607818008941707, 0.8039140004661657, 0.9821306644713012, 0.9148584059763234, 1.319693261450865, 0.9139503744091894, 0.9139503744091894, 0.7423551603023174, 1.175086584622012, 0.9016499087390129, 1.088041181217738, 1.137039275473244, 1.088041181217738, 1.081267050207231, 1.047803715047622, 1.043208066649197),
    @(0.8188218511808876, 0.973014291631896, 0.9873733772220105, 1.061432773325683, 0.8188218511808876, 0.973014291631896, 0.9447071843883361, 1.070354008295701, 0.9781989437905825, 1.051344657150992, 0.8188218511808876, 0.8188218511808876, 0.9873733772220105, 0.9801938344269532, 1.024403075273847, 0.9264031733449313, 1.00727348072653, 0.9264031733449314, 0.9601605733401193, 0.9318928289082729, 0.9856558858732611),
    @(2.88, 0.1, 0.65, 0.66, 2.88, 0.1, 1.29, 0.41, 1.51, 0.24, 2.88, 2.88, 0.65, 0.375, 0.655, 1.21, 0.47, 1.21, 1.0725, 1.434, 0.9675),
    @(0.6899999999999999, 0.21, 1.98, 0.71, 0.6899999999999999, 0.21, 1.47, 1.15, 0.83, 0.47, 0.6899999999999999, 0.6899999999999999, 1.98, 1.095, 1.345, 0.96, 0.9666666666666667, 0.96, 0.8975, 0.8559999999999999, 0.93875),
    @(0.9975350210529828, 0.9958735156987497, 0.9910630543160649, 0.9940494083606265, 0.9975350210529828, 0.9958735156987497, 0.993206838784216, 0.9939777443003533, 0.9957064030097509, 0.9938675598783352, 0.9975220683420907, 0.9975350210529828, 0.9910630543160649, 0.9934682850074072, 0.9925562313383457, 0.9948238636892658, 0.9936619927918137, 0.9948238636892658, 0.994630249857106, 0.9952112040962813, 0.994409943175135)
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item(3 + $r, 3 + $c).Value = $rowVals[$c]
    }
}

